$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.911.37"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +0.03%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.549.99"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.35%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "206.24"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("E9").Value = "  +0.13%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0589"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").Value = "  -0.23%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.770.47"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -0.12%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.549.42"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("E15").Value = "  +0.87%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.907.63"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +0.02%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "61.66"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -0.03%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0706"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +2.89%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "217.28"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  -1.12%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "153.72"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.37%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.95"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -0.16%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +4.78%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.408.64"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("E35").Value = "  +2.91%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.966"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -0.02%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0165"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +0.89%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.530"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("E44").Value = "  +0.89%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "64.54"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("E46").Value = "  +0.42%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.685.07"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -0.09%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "87.46"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("E50").Value = "  +3.67%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0961"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +0.73%  "
